# Add "dummy_capacity" variable to the "Coupling Parameters" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Existing "maximum_investment_capacity_per_year" value changes from 1000000 to 10000
$ws.Range("B13").Value = 10000

# New row 19: dummy_capacity variable
$ws.Range("A19").Value = "dummy_capacity"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = "MW capacity to be assigned to candidate power plants"

# Match the author's final on-screen selection
$ws.Range("B23").Select() | Out-Null
